# Update weekly triaged issues
# Rewrites the component stats table (rows 2-50) with the refreshed
# bug/feature/total counts, re-sorted order, and two new rows that
# extend the data range from A1:D48 to A1:D50.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Menu", 15, 5, 20),
    @("Combobox", 7, 4, 11),
    @("Nav", 8, 1, 9),
    @("DataGrid", 4, 5, 9),
    @("Popover", 6, 2, 8),
    @("Dropdown", 6, 2, 8),
    @("Tree", 6, 0, 6),
    @("Dialog", 4, 1, 5),
    @("Tooltip", 4, 1, 5),
    @("Drawer", 4, 0, 4),
    @("Virtualizer", 4, 0, 4),
    @("Toolbar", 3, 1, 4),
    @("TagPicker", 2, 1, 3),
    @("Table", 3, 0, 3),
    @("TeachingPopover", 3, 0, 3),
    @("Skeleton", 2, 1, 3),
    @("DatePicker", 0, 2, 2),
    @("FluentProvider", 0, 2, 2),
    @("Calendar Compat", 2, 0, 2),
    @("Slider", 0, 2, 2),
    @("MessageBar", 2, 0, 2),
    @("Portal", 2, 0, 2),
    @("Tabs", 2, 0, 2),
    @("Switch", 2, 0, 2),
    @("List", 2, 0, 2),
    @("Popup", 1, 0, 1),
    @("MenuItem", 1, 0, 1),
    @("Label", 1, 0, 1),
    @("AvatarGroup", 1, 0, 1),
    @("Card", 1, 0, 1),
    @("Avatar", 1, 0, 1),
    @("InfoLabel", 1, 0, 1),
    @("Button", 0, 1, 1),
    @("Checkbox", 1, 0, 1),
    @("Image", 0, 1, 1),
    @("Input", 0, 1, 1),
    @("TimePicker", 1, 0, 1),
    @("FocusTrapZone", 1, 0, 1),
    @("Badge", 1, 0, 1),
    @("Accordion", 1, 0, 1),
    @("Toast", 1, 0, 1),
    @("SplitButton", 1, 0, 1),
    @("Coachmark", 0, 0, 0),
    @("Calendar", 0, 0, 0),
    @("Keytip", 0, 0, 0),
    @("Pickers", 0, 0, 0),
    @("Carousel", 0, 0, 0),
    @("Segment", 0, 0, 0),
    @("ColorPicker", 0, 0, 0)
)

$rowCount = $data.Count
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $entry = $data[$i]
    $ws.Cells.Item($r, 1).Value = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]
    $ws.Cells.Item($r, 4).Value = $entry[3]
}
